# Fruta / hortaliza, semanal
# New weekly price report: insert two new rows (new "Perú" origin quality
# grades "1a nueva(o)" / "2a nueva(o)" for the week of 2021-11-09) right
# before the existing "Camote" block, pushing the rest of that block down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 129 — Excel shifts rows 129:134 down to
# 131:136 and copies formatting (incl. the date-cell style on column D) from
# the row above, matching the row-above's number format.
$ws.Rows("129:130").Insert()

# Row 129 — Camote, 1a nueva(o), Perú
$ws.Range("A129").Value = 11
$ws.Range("B129").Value = "Vega Monumental Concepción"
$ws.Range("C129").Value = "Bíobío"
$ws.Range("D129").Value2 = 44509
$ws.Range("E129").Value = 8
$ws.Range("F129").Value = 100112045
$ws.Range("G129").Value = "Zapallo"
$ws.Range("H129").Value = "Camote"
$ws.Range("I129").Value = "1a nueva(o)"
$ws.Range("J129").Value = 400
$ws.Range("K129").Value = 650
$ws.Range("L129").Value = 650
$ws.Range("M129").Value = 650
$ws.Range("N129").Value = "$/kilo (volumen en unidades)"
$ws.Range("O129").Value = "Perú"
$ws.Range("P129").Value = 650
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"

# Row 130 — Camote, 2a nueva(o), Perú
$ws.Range("A130").Value = 11
$ws.Range("B130").Value = "Vega Monumental Concepción"
$ws.Range("C130").Value = "Bíobío"
$ws.Range("D130").Value2 = 44509
$ws.Range("E130").Value = 8
$ws.Range("F130").Value = 100112045
$ws.Range("G130").Value = "Zapallo"
$ws.Range("H130").Value = "Camote"
$ws.Range("I130").Value = "2a nueva(o)"
$ws.Range("J130").Value = 400
$ws.Range("K130").Value = 550
$ws.Range("L130").Value = 550
$ws.Range("M130").Value = 550
$ws.Range("N130").Value = "$/kilo (volumen en unidades)"
$ws.Range("O130").Value = "Perú"
$ws.Range("P130").Value = 550
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = "Hortaliza"
